# Issue #57: Make genre required with PBCore controlled vocabulary.
# Fixture update: add a "Genre" column (T) to the non-manager batch
# ingest manifest spreadsheet, with an "Auction" value for the two
# sample data rows, so the manifest can exercise the genre-override
# behaviour during bibliographic import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (row 2), matching the existing "Genre" header used
# in column D.
$ws.Range("T2").Value = "Genre"

# New per-row genre values (rows 3 and 4) using the PBCore controlled
# vocabulary term "Auction".
$ws.Range("T3").Value = "Auction"
$ws.Range("T4").Value = "Auction"

# Move the active selection to reflect the newly extended used range,
# same as a user who just finished typing into the last new cell.
$ws.Range("T5").Select()
